$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I5").Value = 30
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("J5").Value = 0
$ws.Range("M5").Value = 85
$ws.Range("H5").Value = 30
$ws.Range("K5").Value = 30
$ws.Range("I12").Value = 1000
$ws.Range("M12").Value = -830
$ws.Range("H12").Value = 2333.3333
$ws.Range("K12").Value = 1000
$ws.Range("I15").Value = 2299.175
$ws.Range("M15").Value = -6728.525000000001
$ws.Range("H15").Value = 2299.175
$ws.Range("K15").Value = 6897.525000000001
$ws.Range("L16").Value = 20000
$ws.Range("N16").Value = -20460
$ws.Range("J16").Value = 20000
$ws.Range("H16").Value = 20000
$ws.Range("L17").Value = 5472.9
$ws.Range("N17").Value = -5808.9
$ws.Range("J17").Value = 1824.3
$ws.Range("H17").Value = 2108.1667
$ws.Range("L40").Value = 10004
$ws.Range("N40").Value = -10354
$ws.Range("J40").Value = 10004
$ws.Range("H40").Value = 7644
$ws.Range("I46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("I62").Value = 2339.1667
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
$ws.Range("J62").Value = 5000
$ws.Range("M62").Value = -1715.1667
$ws.Range("H62").Value = 2719.2856
$ws.Range("K62").Value = 2339.1667
$ws.Range("I65").Value = 2339.1667
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
$ws.Range("J65").Value = 5000
$ws.Range("M65").Value = -8575.833500000001
$ws.Range("H65").Value = 2719.2856
$ws.Range("K65").Value = 11695.8335
$ws.Range("I70").Value = 3399.7
$ws.Range("L70").Value = 15814.2861
$ws.Range("N70").Value = -16354.2861
$ws.Range("J70").Value = 5271.4287
$ws.Range("M70").Value = -9929.099999999999
$ws.Range("H70").Value = 4170.4116
$ws.Range("K70").Value = 10199.1
$ws.Range("I73").Value = 3399.7
$ws.Range("L73").Value = 15814.2861
$ws.Range("N73").Value = -17686.2861
$ws.Range("J73").Value = 5271.4287
$ws.Range("M73").Value = -9263.099999999999
$ws.Range("H73").Value = 4170.4116
$ws.Range("K73").Value = 10199.1
$ws.Range("I76").Value = 4840.6
$ws.Range("M76").Value = -4525.6
$ws.Range("H76").Value = 5033.6665
$ws.Range("K76").Value = 4840.6
$ws.Range("I79").Value = 4840.6
$ws.Range("M79").Value = -3748.6
$ws.Range("H79").Value = 5033.6665
$ws.Range("K79").Value = 4840.6
$ws.Range("I80").Value = 496.5
$ws.Range("M80").Value = -491.5
$ws.Range("H80").Value = 496.5
$ws.Range("K80").Value = 1489.5
$ws.Range("I83").Value = 496.5
$ws.Range("M83").Value = 523.5
$ws.Range("H83").Value = 496.5
$ws.Range("K83").Value = 4468.5
$ws.Range("I86").Value = 111117360
$ws.Range("M86").Value = -111116237
$ws.Range("H86").Value = 76930136
$ws.Range("K86").Value = 111117360
$ws.Range("L87").Value = 179987.5
$ws.Range("N87").Value = -182483.5
$ws.Range("J87").Value = 179987.5
$ws.Range("H87").Value = 179987.5
$ws.Range("L88").Value = 10771.909
$ws.Range("N88").Value = -11583.909
$ws.Range("J88").Value = 10771.909
$ws.Range("H88").Value = 5271832.5
$ws.Range("I89").Value = 111117360
$ws.Range("M89").Value = -555581184
$ws.Range("H89").Value = 76930136
$ws.Range("K89").Value = 555586800
$ws.Range("L90").Value = 539962.5
$ws.Range("N90").Value = -552442.5
$ws.Range("J90").Value = 179987.5
$ws.Range("H90").Value = 179987.5
$ws.Range("L91").Value = 10771.909
$ws.Range("N91").Value = -13579.909
$ws.Range("J91").Value = 10771.909
$ws.Range("H91").Value = 5271832.5
$ws.Range("I98").Value = 3903.2104
$ws.Range("L98").Value = 6302
$ws.Range("N98").Value = -9298
$ws.Range("J98").Value = 6302
$ws.Range("M98").Value = -2405.2104
$ws.Range("H98").Value = 4230.3184
$ws.Range("K98").Value = 3903.2104
$ws.Range("I106").Value = 2005
$ws.Range("L106").Value = 3750
$ws.Range("N106").Value = -5012
$ws.Range("J106").Value = 3750
$ws.Range("M106").Value = -1374
$ws.Range("H106").Value = 3168.3333
$ws.Range("K106").Value = 2005
$ws.Range("L121").Value = 10700.4999
$ws.Range("N121").Value = -14194.4999
$ws.Range("J121").Value = 3566.8333
$ws.Range("H121").Value = 3566.8333
$ws.Range("I122").Value = 3903.2104
$ws.Range("L122").Value = 18906
$ws.Range("N122").Value = -23806
$ws.Range("J122").Value = 6302
$ws.Range("M122").Value = -9259.6312
$ws.Range("H122").Value = 4230.3184
$ws.Range("K122").Value = 11709.6312
$ws.Range("I127").Value = 1009.05554
$ws.Range("L127").Value = 5999.6667
$ws.Range("N127").Value = -15919.6667
$ws.Range("J127").Value = 1999.8889
$ws.Range("M127").Value = 1932.83338
$ws.Range("H127").Value = 1339.3334
$ws.Range("K127").Value = 3027.16662
$ws.Range("I131").Value = 1263382.8
$ws.Range("L131").Value = 34657.5
$ws.Range("N131").Value = -44737.5
$ws.Range("J131").Value = 11552.5
$ws.Range("M131").Value = -3785108.4
$ws.Range("H131").Value = 1013016.7
$ws.Range("K131").Value = 3790148.4
$ws.Range("I132").Value = 8930.056
$ws.Range("M132").Value = -24260.168
$ws.Range("H132").Value = 4686822.5
$ws.Range("K132").Value = 26790.168
$ws.Range("I135").Value = 1312.84
$ws.Range("L135").Value = 13580.1
$ws.Range("N135").Value = -18650.1
$ws.Range("J135").Value = 1508.9
$ws.Range("M135").Value = -9280.559999999999
$ws.Range("H135").Value = 1368.8572
$ws.Range("K135").Value = 11815.56
$ws.Range("I137").Value = 2085223.2
$ws.Range("L137").Value = 14810.1
$ws.Range("N137").Value = -19910.1
$ws.Range("J137").Value = 4936.7
$ws.Range("M137").Value = -6253119.6
$ws.Range("H137").Value = 1473374.2
$ws.Range("K137").Value = 6255669.6
$ws.Range("I138").Value = 1170.6296
$ws.Range("L138").Value = 9418.068600000001
$ws.Range("N138").Value = -19698.0686
$ws.Range("J138").Value = 3139.3562
$ws.Range("M138").Value = 1628.1112
$ws.Range("H138").Value = 2607.8
$ws.Range("K138").Value = 3511.8888
$ws.Range("I141").Value = 1456.9333
$ws.Range("L141").Value = 6297
$ws.Range("N141").Value = -16657
$ws.Range("J141").Value = 2099
$ws.Range("M141").Value = 809.2001
$ws.Range("H141").Value = 1497.0625
$ws.Range("K141").Value = 4370.7999

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I22").Value = 6999
$ws.Range("M22").Value = -6700
$ws.Range("H22").Value = 6999
$ws.Range("K22").Value = 6999
$ws.Range("I32").Value = 2370.5688
$ws.Range("M32").Value = -2083.5688
$ws.Range("H32").Value = 2378.0508
$ws.Range("K32").Value = 2370.5688
$ws.Range("I36").Value = 10012.5
$ws.Range("M36").Value = -9666.5
$ws.Range("H36").Value = 10012.5
$ws.Range("K36").Value = 10012.5
$ws.Range("I45").Value = 48535.445
$ws.Range("M45").Value = -48158.445
$ws.Range("H45").Value = 48535.445
$ws.Range("K45").Value = 48535.445
$ws.Range("I61").Value = 1793.8334
$ws.Range("L61").Value = 4039
$ws.Range("N61").Value = -4463
$ws.Range("J61").Value = 4039
$ws.Range("M61").Value = -1581.8334
$ws.Range("H61").Value = 2916.4167
$ws.Range("K61").Value = 1793.8334
$ws.Range("I74").Value = 207782.78
$ws.Range("L74").Value = 1904.5
$ws.Range("N74").Value = -3652.5
$ws.Range("J74").Value = 1904.5
$ws.Range("M74").Value = -206908.78
$ws.Range("H74").Value = 193584.28
$ws.Range("K74").Value = 207782.78
$ws.Range("I77").Value = 207782.78
$ws.Range("L77").Value = 9522.5
$ws.Range("N77").Value = -18258.5
$ws.Range("J77").Value = 1904.5
$ws.Range("M77").Value = -1034545.9
$ws.Range("H77").Value = 193584.28
$ws.Range("K77").Value = 1038913.9
$ws.Range("I88").Value = 2171.7273
$ws.Range("L88").Value = 3205.15
$ws.Range("N88").Value = -4017.15
$ws.Range("J88").Value = 3205.15
$ws.Range("M88").Value = -1765.7273
$ws.Range("H88").Value = 2838.4517
$ws.Range("K88").Value = 2171.7273
$ws.Range("I91").Value = 2171.7273
$ws.Range("L91").Value = 3205.15
$ws.Range("N91").Value = -6013.15
$ws.Range("J91").Value = 3205.15
$ws.Range("M91").Value = -767.7273
$ws.Range("H91").Value = 2838.4517
$ws.Range("K91").Value = 2171.7273
$ws.Range("I132").Value = 1974088.1
$ws.Range("M132").Value = -5919734.300000001
$ws.Range("H132").Value = 9119008
$ws.Range("K132").Value = 5922264.300000001
$ws.Range("I136").Value = 1793.8334
$ws.Range("L136").Value = 12117
$ws.Range("N136").Value = -17217
$ws.Range("J136").Value = 4039
$ws.Range("M136").Value = -2831.5002
$ws.Range("H136").Value = 2916.4167
$ws.Range("K136").Value = 5381.5002
$ws.Range("L139").Value = 75588.8
$ws.Range("N139").Value = -85868.8
$ws.Range("J139").Value = 75588.8
$ws.Range("H139").Value = 75588.8
$ws.Range("L140").Value = 72599.664
$ws.Range("N140").Value = -82959.664
$ws.Range("J140").Value = 72599.664
$ws.Range("H140").Value = 72599.664

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("L80").Value = 732
$ws.Range("N80").Value = -2728
$ws.Range("J80").Value = 732
$ws.Range("H80").Value = 692
$ws.Range("L83").Value = 3660
$ws.Range("N83").Value = -13644
$ws.Range("J83").Value = 732
$ws.Range("H83").Value = 692
$ws.Range("I99").Value = 102812.3
$ws.Range("M99").Value = -101314.3
$ws.Range("H99").Value = 62847.59
$ws.Range("K99").Value = 102812.3
$ws.Range("I134").Value = 3082.88
$ws.Range("L134").Value = 11657.0001
$ws.Range("N134").Value = -16727.0001
$ws.Range("J134").Value = 3885.6667
$ws.Range("M134").Value = -6713.639999999999
$ws.Range("H134").Value = 3238.258
$ws.Range("K134").Value = 9248.639999999999

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I22").Value = 2159.8
$ws.Range("M22").Value = -1809.8
$ws.Range("H22").Value = 2042.5714
$ws.Range("K22").Value = 2159.8
$ws.Range("I31").Value = 4395.5
$ws.Range("L31").Value = 6792.3076
$ws.Range("N31").Value = -7382.3076
$ws.Range("J31").Value = 6792.3076
$ws.Range("M31").Value = -4100.5
$ws.Range("H31").Value = 5750.2173
$ws.Range("K31").Value = 4395.5
$ws.Range("I34").Value = 4395.5
$ws.Range("L34").Value = 6792.3076
$ws.Range("N34").Value = -7196.3076
$ws.Range("J34").Value = 6792.3076
$ws.Range("M34").Value = -4193.5
$ws.Range("H34").Value = 5750.2173
$ws.Range("K34").Value = 4395.5
$ws.Range("L37").Value = 5000
$ws.Range("N37").Value = -5214
$ws.Range("J37").Value = 5000
$ws.Range("H37").Value = 5000
$ws.Range("I62").Value = 20003100
$ws.Range("M62").Value = -20002476
$ws.Range("H62").Value = 11114611
$ws.Range("K62").Value = 20003100
$ws.Range("I65").Value = 20003100
$ws.Range("M65").Value = -100012380
$ws.Range("H65").Value = 11114611
$ws.Range("K65").Value = 100015500
$ws.Range("I86").Value = 10000
$ws.Range("L86").Value = 9900
$ws.Range("N86").Value = -12146
$ws.Range("J86").Value = 9900
$ws.Range("M86").Value = -8877
$ws.Range("H86").Value = 9950
$ws.Range("K86").Value = 10000
$ws.Range("I89").Value = 10000
$ws.Range("L89").Value = 49500
$ws.Range("N89").Value = -60732
$ws.Range("J89").Value = 9900
$ws.Range("M89").Value = -44384
$ws.Range("H89").Value = 9950
$ws.Range("K89").Value = 50000
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("J99").Value = 0
$ws.Range("H99").Value = 70000
$ws.Range("L100").Value = 100000
$ws.Range("N100").Value = -102164
$ws.Range("J100").Value = 100000
$ws.Range("H100").Value = 100000
$ws.Range("L107").Value = 1200.2858
$ws.Range("N107").Value = -5040.2858
$ws.Range("J107").Value = 1200.2858
$ws.Range("H107").Value = 3125830.5
$ws.Range("L109").Value = 48984.5
$ws.Range("N109").Value = -51064.5
$ws.Range("J109").Value = 48984.5
$ws.Range("H109").Value = 48984.5
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("J126").Value = 0
$ws.Range("H126").Value = 70000
$ws.Range("I132").Value = 1669.2354
$ws.Range("L132").Value = 100011732
$ws.Range("N132").Value = -100016792
$ws.Range("J132").Value = 33337244
$ws.Range("M132").Value = -2477.706200000001
$ws.Range("H132").Value = 12348178
$ws.Range("K132").Value = 5007.706200000001
$ws.Range("I134").Value = 2464.1853
$ws.Range("L134").Value = 11399.4
$ws.Range("N134").Value = -16469.4
$ws.Range("J134").Value = 3799.8
$ws.Range("M134").Value = -4857.5559
$ws.Range("H134").Value = 2672.875
$ws.Range("K134").Value = 7392.5559

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I5").Value = 628.85
$ws.Range("L5").Value = 2268.9999
$ws.Range("N5").Value = -2492.9999
$ws.Range("J5").Value = 756.3333
$ws.Range("M5").Value = -1774.55
$ws.Range("H5").Value = 645.4783
$ws.Range("K5").Value = 1886.55
$ws.Range("I17").Value = 264
$ws.Range("M17").Value = -623
$ws.Range("H17").Value = 264
$ws.Range("K17").Value = 792
$ws.Range("I34").Value = 282.33334
$ws.Range("L34").Value = 11149.3998
$ws.Range("N34").Value = -11317.3998
$ws.Range("J34").Value = 3716.4666
$ws.Range("M34").Value = -763.0000200000001
$ws.Range("H34").Value = 3144.111
$ws.Range("K34").Value = 847.0000200000001
$ws.Range("L39").Value = 28644.6
$ws.Range("N39").Value = -29232.6
$ws.Range("J39").Value = 9548.200000000001
$ws.Range("H39").Value = 8355.857
$ws.Range("I60").Value = 1429481.8
$ws.Range("L60").Value = 13113.819
$ws.Range("N60").Value = -13615.819
$ws.Range("J60").Value = 4371.273
$ws.Range("M60").Value = -4288194.4
$ws.Range("H60").Value = 802433.2
$ws.Range("K60").Value = 4288445.4
$ws.Range("L98").Value = 1811.25
$ws.Range("N98").Value = -4807.25
$ws.Range("J98").Value = 603.75
$ws.Range("H98").Value = 512.8
$ws.Range("L121").Value = 2989.5
$ws.Range("N121").Value = -5609.5
$ws.Range("J121").Value = 996.5
$ws.Range("H121").Value = 20400558
$ws.Range("L122").Value = 6927.3
$ws.Range("N122").Value = -11827.3
$ws.Range("J122").Value = 769.7
$ws.Range("H122").Value = 502.3846
$ws.Range("I131").Value = 62785.8
$ws.Range("L131").Value = 5876.833500000001
$ws.Range("N131").Value = -15956.8335
$ws.Range("J131").Value = 1958.9445
$ws.Range("M131").Value = -183317.4
$ws.Range("H131").Value = 15182.174
$ws.Range("K131").Value = 188357.4
$ws.Range("I135").Value = 628.85
$ws.Range("L135").Value = 6806.9997
$ws.Range("N135").Value = -11876.9997
$ws.Range("J135").Value = 756.3333
$ws.Range("M135").Value = -3124.650000000001
$ws.Range("H135").Value = 645.4783
$ws.Range("K135").Value = 5659.650000000001
$ws.Range("L137").Value = 8906.499899999999
$ws.Range("N137").Value = -19106.4999
$ws.Range("J137").Value = 2968.8333
$ws.Range("H137").Value = 1963.7037

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I19").Value = 500
$ws.Range("M19").Value = -212
$ws.Range("H19").Value = 500
$ws.Range("K19").Value = 500
$ws.Range("I46").Value = 39500
$ws.Range("M46").Value = -39344
$ws.Range("H46").Value = 33729.2
$ws.Range("K46").Value = 39500
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("J52").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("I70").Value = 27892216
$ws.Range("M70").Value = -27891946
$ws.Range("H70").Value = 20088172
$ws.Range("K70").Value = 27892216
$ws.Range("I73").Value = 27892216
$ws.Range("M73").Value = -27891280
$ws.Range("H73").Value = 20088172
$ws.Range("K73").Value = 27892216
$ws.Range("I122").Value = 2342.9092
$ws.Range("L122").Value = 23452.3329
$ws.Range("N122").Value = -28352.3329
$ws.Range("J122").Value = 7817.4443
$ws.Range("M122").Value = -4578.7276
$ws.Range("H122").Value = 3932.2903
$ws.Range("K122").Value = 7028.7276
$ws.Range("I132").Value = 3875.5
$ws.Range("L132").Value = 11400
$ws.Range("N132").Value = -16460
$ws.Range("J132").Value = 3800
$ws.Range("M132").Value = -9096.5
$ws.Range("H132").Value = 3860.4
$ws.Range("K132").Value = 11626.5
$ws.Range("L135").Value = 69999
$ws.Range("N135").Value = -80139
$ws.Range("J135").Value = 69999
$ws.Range("H135").Value = 69999

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I16").Value = 827.6
$ws.Range("M16").Value = -657.6
$ws.Range("H16").Value = 1723.7273
$ws.Range("K16").Value = 827.6
$ws.Range("I30").Value = 1852.6666
$ws.Range("M30").Value = -1744.6666
$ws.Range("H30").Value = 1852.6666
$ws.Range("K30").Value = 1852.6666
$ws.Range("I46").Value = 2079.3333
$ws.Range("L46").Value = 2500
$ws.Range("N46").Value = -2876
$ws.Range("J46").Value = 2500
$ws.Range("M46").Value = -1891.3333
$ws.Range("H46").Value = 2219.5557
$ws.Range("K46").Value = 2079.3333
$ws.Range("I61").Value = 3948.8333
$ws.Range("M61").Value = -3746.8333
$ws.Range("H61").Value = 14658.8
$ws.Range("K61").Value = 3948.8333
$ws.Range("I93").Value = 2725.4614
$ws.Range("M93").Value = -1477.4614
$ws.Range("H93").Value = 2789.1177
$ws.Range("K93").Value = 2725.4614
$ws.Range("I100").Value = 2447.5
$ws.Range("M100").Value = -1906.5
$ws.Range("H100").Value = 2798.6
$ws.Range("K100").Value = 2447.5
$ws.Range("I113").Value = 3948.8333
$ws.Range("M113").Value = -1778.8333
$ws.Range("H113").Value = 14658.8
$ws.Range("K113").Value = 3948.8333
$ws.Range("I122").Value = 4249.25
$ws.Range("M122").Value = -10297.75
$ws.Range("H122").Value = 3499.3333
$ws.Range("K122").Value = 12747.75
$ws.Range("I132").Value = 2380.3684
$ws.Range("L132").Value = 11928.6921
$ws.Range("N132").Value = -16988.6921
$ws.Range("J132").Value = 3976.2307
$ws.Range("M132").Value = -4611.1052
$ws.Range("H132").Value = 3028.6875
$ws.Range("K132").Value = 7141.1052
$ws.Range("L133").Value = 103998.336
$ws.Range("N133").Value = -109058.336
$ws.Range("J133").Value = 103998.336
$ws.Range("H133").Value = 103998.336
$ws.Range("I136").Value = 2705.1765
$ws.Range("L136").Value = 10212.4284
$ws.Range("N136").Value = -15312.4284
$ws.Range("J136").Value = 3404.1428
$ws.Range("M136").Value = -5565.529500000001
$ws.Range("H136").Value = 2909.0417
$ws.Range("K136").Value = 8115.529500000001

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I14").Value = 0
$ws.Range("L14").Value = 30004.4
$ws.Range("N14").Value = -30340.4
$ws.Range("J14").Value = 30004.4
$ws.Range("M14").ClearContents()
$ws.Range("K14").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("J24").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("I27").Value = 49569
$ws.Range("L27").Value = 49998
$ws.Range("N27").Value = -50136
$ws.Range("J27").Value = 49998
$ws.Range("M27").Value = -49500
$ws.Range("H27").Value = 49783.5
$ws.Range("K27").Value = 49569
$ws.Range("L45").Value = 19382.727
$ws.Range("N45").Value = -20364.727
$ws.Range("J45").Value = 19382.727
$ws.Range("H45").Value = 16198.5
$ws.Range("L46").Value = 49797.5
$ws.Range("N46").Value = -50259.5
$ws.Range("J46").Value = 49797.5
$ws.Range("H46").Value = 49797.5
$ws.Range("I58").Value = 3997.6667
$ws.Range("M58").Value = -3689.6667
$ws.Range("H58").Value = 9248.25
$ws.Range("K58").Value = 3997.6667
$ws.Range("I62").Value = 4998
$ws.Range("L62").Value = 7999
$ws.Range("N62").Value = -9247
$ws.Range("J62").Value = 7999
$ws.Range("M62").Value = -4374
$ws.Range("H62").Value = 7248.75
$ws.Range("K62").Value = 4998
$ws.Range("I65").Value = 4998
$ws.Range("L65").Value = 39995
$ws.Range("N65").Value = -46235
$ws.Range("J65").Value = 7999
$ws.Range("M65").Value = -21870
$ws.Range("H65").Value = 7248.75
$ws.Range("K65").Value = 24990
$ws.Range("I81").Value = 1969.7142
$ws.Range("M81").Value = -2878.4284
$ws.Range("H81").Value = 4098.5454
$ws.Range("K81").Value = 3939.4284
$ws.Range("I84").Value = 1969.7142
$ws.Range("M84").Value = -14393.142
$ws.Range("H84").Value = 4098.5454
$ws.Range("K84").Value = 19697.142
$ws.Range("I107").Value = 536.5789
$ws.Range("M107").Value = 310.2633000000001
$ws.Range("H107").Value = 594.8387
$ws.Range("K107").Value = 1609.7367
$ws.Range("I113").Value = 637.7143
$ws.Range("M113").Value = 256.8571000000002
$ws.Range("H113").Value = 761.3929000000001
$ws.Range("K113").Value = 1913.1429
$ws.Range("I115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("H115").Value = 49999
$ws.Range("K115").Value = 0
$ws.Range("I122").Value = 6165
$ws.Range("M122").Value = -16045
$ws.Range("H122").Value = 19236152
$ws.Range("K122").Value = 18495
$ws.Range("I132").Value = 1516.3726
$ws.Range("L132").Value = 11369.0001
$ws.Range("N132").Value = -16429.0001
$ws.Range("J132").Value = 3789.6667
$ws.Range("M132").Value = -2019.1178
$ws.Range("H132").Value = 1857.3667
$ws.Range("K132").Value = 4549.1178
$ws.Range("L133").Value = 83861.625
$ws.Range("N133").Value = -93981.625
$ws.Range("J133").Value = 83861.625
$ws.Range("H133").Value = 83861.625
$ws.Range("L134").Value = 149392.5
$ws.Range("N134").Value = -154462.5
$ws.Range("J134").Value = 49797.5
$ws.Range("H134").Value = 49797.5
